$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Program to cyclically rotate an array by one"
$ws.Range("B13").Value = "CyclicallyRotate"

$ws.Range("B13").Select()
